$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) — row 2
$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 34.35779879069297

# Regenerated s_vals data (filtered save games) — row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 15.88780690183548
